$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38:C38").NumberFormat = "@"
$ws.Range("A38").Value = "2025-09-24"
$ws.Range("B38").Value = "21:20:39"
$ws.Range("C38").Value = "1.00 EUR = 1,624.6901"
